$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
$notes.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Speaker notes"
